$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Улучшение"
$ws.Range("B19").Value = "Добавить больше теоретических материалов (Как создавать ассоциации)"

$ws.Range("B20").Select()
